$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMP")

# Map of row -> column -> new value, for the margin rows updated in this commit.
$updates = @{
    14 = @{ "D" = 0.9933; "E" = 0.9948; "F" = 0.9984; "G" = 1.0 }
    15 = @{ "D" = 0.2502; "E" = 0.3086; "F" = 0.3788; "G" = 0.2402 }
    16 = @{ "D" = 0.18;   "E" = 0.2423; "F" = 0.3196; "G" = 0.1721 }
    17 = @{ "D" = 0.1506; "E" = 0.2022; "F" = 0.2746; "G" = 0.146 }
    18 = @{ "D" = 0.3961; "E" = 0.4708; "F" = 0.5613; "G" = 0.1695 }
    29 = @{ "D" = 0.2669; "E" = 0.3238; "F" = 0.3931; "G" = 0.2543 }
    30 = @{ "D" = 0.4085; "E" = 0.4819; "F" = 0.5722; "G" = 0.1805 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
